$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.151.21'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '3.740.80'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '592.57'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').Value = '166.99'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('D7').Value = '3.737.63'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').Value = '6.40'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '4.366.67'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '3.727.66'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '68.083.46'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('D19').Value = '6.99'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '464.79'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '0.695'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '83.76'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000147'
$ws.Range('E25').Value = '  +8.53%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').Value = '10.02'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '2.77'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').Value = '29.76'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('E33').Value = '  -3.17%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  --%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '9.09'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('D36').Value = '3.694.31'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').Value = '3.44'
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').Value = '0.992'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D44').Value = '44.02'
$ws.Range('E44').Value = '  +15.73%  '
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').Value = '46.74'
$ws.Range('E46').Value = '  +3.25%  '
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').Value = '388.85'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = '143.91'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').Value = '2.749.59'
$ws.Range('E51').Value = '  +2.62%  '

Write-Output "Applied 86 cell updates"
